$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 16377.333
$ws.Range("I21").Value = 5000
$ws.Range("J21").Value = 17799.5
$ws.Range("K21").Value = 5000
$ws.Range("L21").Value = 17799.5
$ws.Range("M21").Value = -4532
$ws.Range("N21").Value = -18735.5
$ws.Range("H23").Value = 16377.333
$ws.Range("I23").Value = 5000
$ws.Range("J23").Value = 17799.5
$ws.Range("K23").Value = 5000
$ws.Range("L23").Value = 17799.5
$ws.Range("M23").Value = -4766
$ws.Range("N23").Value = -18267.5
$ws.Range("H28").Value = 11216.5
$ws.Range("I28").Value = 1469.75
$ws.Range("J28").Value = 20963.25
$ws.Range("K28").Value = 1469.75
$ws.Range("L28").Value = 20963.25
$ws.Range("M28").Value = -984.75
$ws.Range("N28").Value = -21933.25
$ws.Range("H33").Value = 671.2727
$ws.Range("I33").Value = 861.875
$ws.Range("J33").Value = 163
$ws.Range("K33").Value = 861.875
$ws.Range("L33").Value = 163
$ws.Range("M33").Value = -632.875
$ws.Range("N33").Value = -621
$ws.Range("H55").Value = 1047.8462
$ws.Range("I55").Value = 1129.091
$ws.Range("J55").Value = 601
$ws.Range("K55").Value = 1129.091
$ws.Range("L55").Value = 601
$ws.Range("M55").Value = -915.0909999999999
$ws.Range("N55").Value = -1029
$ws.Range("H74").Value = 4182.364
$ws.Range("I74").Value = 4001
$ws.Range("J74").Value = 4400
$ws.Range("K74").Value = 4001
$ws.Range("L74").Value = 4400
$ws.Range("M74").Value = -3065
$ws.Range("N74").Value = -6272
$ws.Range("H77").Value = 4182.364
$ws.Range("I77").Value = 4001
$ws.Range("J77").Value = 4400
$ws.Range("K77").Value = 20005
$ws.Range("L77").Value = 22000
$ws.Range("M77").Value = -15325
$ws.Range("N77").Value = -31360
$ws.Range("H113").Value = 2875
$ws.Range("J113").Value = 3200
$ws.Range("L113").Value = 3200
$ws.Range("N113").Value = -9708
$ws.Range("H116").Value = 3133.2
$ws.Range("I116").Value = 2857.1428
$ws.Range("J116").Value = 3374.75
$ws.Range("K116").Value = 2857.1428
$ws.Range("L116").Value = 3374.75
$ws.Range("M116").Value = 584.8571999999999
$ws.Range("N116").Value = -10258.75
$ws.Range("H129").Value = 1204.0286
$ws.Range("J129").Value = 1789.1904
$ws.Range("L129").Value = 5367.5712
$ws.Range("N129").Value = -15367.5712
$ws.Range("H132").Value = 2465.5938
$ws.Range("I132").Value = 2253.3572
$ws.Range("J132").Value = 3951.25
$ws.Range("K132").Value = 6760.071599999999
$ws.Range("L132").Value = 11853.75
$ws.Range("M132").Value = -4230.071599999999
$ws.Range("N132").Value = -16913.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1377.75
$ws.Range("I2").Value = 1405.5
$ws.Range("K2").Value = 1405.5
$ws.Range("M2").Value = -1292.5
$ws.Range("H32").Value = 1192208.5
$ws.Range("I32").Value = 1247863.5
$ws.Range("K32").Value = 1247863.5
$ws.Range("M32").Value = -1247576.5
$ws.Range("H63").Value = 156711.5
$ws.Range("I63").Value = 461000
$ws.Range("J63").Value = 4567.25
$ws.Range("K63").Value = 461000
$ws.Range("L63").Value = 4567.25
$ws.Range("M63").Value = -460314
$ws.Range("N63").Value = -5939.25
$ws.Range("H66").Value = 156711.5
$ws.Range("I66").Value = 461000
$ws.Range("J66").Value = 4567.25
$ws.Range("K66").Value = 2305000
$ws.Range("L66").Value = 22836.25
$ws.Range("M66").Value = -2301568
$ws.Range("N66").Value = -29700.25
$ws.Range("H116").Value = 1377.75
$ws.Range("I116").Value = 1405.5
$ws.Range("K116").Value = 1405.5
$ws.Range("M116").Value = 888.5
$ws.Range("H132").Value = 15400617
$ws.Range("I132").Value = 16670.666
$ws.Range("J132").Value = 38476536
$ws.Range("K132").Value = 50011.99800000001
$ws.Range("L132").Value = 115429608
$ws.Range("M132").Value = -47481.99800000001
$ws.Range("N132").Value = -115434668
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1377.75
$ws.Range("I3").Value = 1405.5
$ws.Range("K3").Value = 1405.5
$ws.Range("M3").Value = -1291.5
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2556.6667
$ws.Range("I16").Value = 2011
$ws.Range("J16").Value = 2624.875
$ws.Range("K16").Value = 2011
$ws.Range("L16").Value = 2624.875
$ws.Range("M16").Value = -1724
$ws.Range("N16").Value = -3198.875
$ws.Range("H31").Value = 5167.5537
$ws.Range("I31").Value = 1372.2
$ws.Range("J31").Value = 7276.0835
$ws.Range("K31").Value = 1372.2
$ws.Range("L31").Value = 7276.0835
$ws.Range("M31").Value = -1077.2
$ws.Range("N31").Value = -7866.0835
$ws.Range("H34").Value = 5167.5537
$ws.Range("I34").Value = 1372.2
$ws.Range("J34").Value = 7276.0835
$ws.Range("K34").Value = 1372.2
$ws.Range("L34").Value = 7276.0835
$ws.Range("M34").Value = -1170.2
$ws.Range("N34").Value = -7680.0835
$ws.Range("H113").Value = 2556.6667
$ws.Range("I113").Value = 2011
$ws.Range("J113").Value = 2624.875
$ws.Range("K113").Value = 2011
$ws.Range("L113").Value = 2624.875
$ws.Range("M113").Value = 159
$ws.Range("N113").Value = -6964.875
$ws.Range("H129").Value = 49999
$ws.Range("J129").Value = 49999
$ws.Range("L129").Value = 49999
$ws.Range("N129").Value = -59999
$ws.Range("H132").Value = 12823601
$ws.Range("I132").Value = 3206
$ws.Range("J132").Value = 33336234
$ws.Range("K132").Value = 9618
$ws.Range("L132").Value = 100008702
$ws.Range("M132").Value = -7088
$ws.Range("N132").Value = -100013762
$ws.Range("H141").Value = 118713.234
$ws.Range("J141").Value = 114648.9
$ws.Range("L141").Value = 114648.9
$ws.Range("N141").Value = -125008.9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 8688.625
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 8688.625
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 26065.875
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -28311.875
$ws.Range("H84").Value = 8688.625
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 8688.625
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 78197.625
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -89429.625
$ws.Range("H92").Value = 850
$ws.Range("I92").Value = 850
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 2550
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -1302
$ws.Range("N92").ClearContents()
$ws.Range("H122").Value = 3077.195
$ws.Range("I122").Value = 368.53845
$ws.Range("K122").Value = 3316.84605
$ws.Range("M122").Value = -866.8460500000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 10007
$ws.Range("I132").Value = 10000
$ws.Range("K132").Value = 30000
$ws.Range("M132").Value = -27470
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").ClearContents()
$ws.Range("H40").Value = 2239.6667
$ws.Range("I40").Value = 1862.7273
$ws.Range("J40").Value = 3276.25
$ws.Range("K40").Value = 1862.7273
$ws.Range("L40").Value = 3276.25
$ws.Range("M40").Value = -1726.7273
$ws.Range("N40").Value = -3548.25
$ws.Range("H42").Value = 4800000
$ws.Range("J42").Value = 0
$ws.Range("L42").Value = 0
$ws.Range("N42").ClearContents()
$ws.Range("H46").Value = 888.6667
$ws.Range("I46").Value = 888.6667
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 888.6667
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -700.6667
$ws.Range("N46").ClearContents()
$ws.Range("H49").Value = 4800000
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H55").Value = 281.64
$ws.Range("J55").Value = 352.41666
$ws.Range("L55").Value = 352.41666
$ws.Range("N55").Value = -698.41666
$ws.Range("H122").Value = 3597.257
$ws.Range("J122").Value = 3898.182
$ws.Range("L122").Value = 11694.546
$ws.Range("N122").Value = -16594.546
